$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "Are there any required standards in effect, implementation language, policies for database integrity, resource limits, operating environment(s) and so on?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Are there any required standards in effect, implementation language, policies for database integrity, resource limits, operating environment(s) and so on",
    2)

Write-Host "Find/Replace result: $found"
